# Regenerate save_data to use K (strikeouts) instead of Strike# (swinging strikes count),
# recalculated std/mean, and write the new s_vals into column G (K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-33 (replacing the old Strike# values)
$newK = @(3,3,8,0,4,3,3,5,6,5,7,12,7,5,7,6,6,6,8,5,3,8,11,7,7,1,3,6,13,6,0,3)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $newK[$i]
}
